# UsersPostAPI.xlsx - "fixed excel and modified schema path"
#
# Semantic changes (per canonical-OOXML diff):
#  - Sheet "PostUsersExistOneField":
#      A2: Anshul       -> Bauwa
#      L2: PT3647       -> PT9815
#      D3: 1234567039   -> 1234567079
#      L3: PT6446       -> PT3732
#      E4: abcn@xyz.com -> abcr@xyz.com
#      L4: PT3270       -> PT1214
#    and this sheet becomes the active/selected tab (selection stays E4).
#  - Sheet "PostUsers":
#      A2: Kusharg -> Muskan
#      L2: PT9732  -> (cleared)
#      stray row 19 (C19 = "s") is removed entirely
#      this sheet is no longer the active tab; its selection moves to A2
#
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("PostUsersExistOneField")
$ws2 = $wb.Worksheets.Item("PostUsers")

# --- PostUsersExistOneField edits ---
$ws1.Range("A2").Value = "Bauwa"
$ws1.Range("L2").Value = "PT9815"
$ws1.Range("D3").Value = 1234567079
$ws1.Range("L3").Value = "PT3732"
$ws1.Range("E4").Value = "abcr@xyz.com"
$ws1.Range("L4").Value = "PT1214"

# --- PostUsers edits ---
$ws2.Range("A2").Value = "Muskan"
$ws2.Range("L2").ClearContents()

# Remove the stray leftover row (C19 = "s") and shrink the used range back to A1:L2
$ws2.Rows.Item(19).Delete()

# Selections: PostUsers' selection moves to A2 and it is no longer the active tab;
# PostUsersExistOneField becomes the active tab (selection stays at E4).
$ws2.Range("A2").Select()
$ws1.Activate()
$ws1.Range("E4").Select()
